$wb = $excel.ActiveWorkbook

# The 02a80f9e-... file has been handed back. Update its status from
# "Ready for handoff" to "Handed back: in sync with en-US" on every sheet,
# and record the new "Latest Handback DateTime" on the per-locale sheets.

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $statusText
$zhcn.Range("G3").Value = "2016-01-15 07:55:25"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $statusText
$dede.Range("G3").Value = "2016-01-15 07:55:42"
